# Re-colour the deck's theme (ppt/theme/theme1.xml, the theme actually used
# by the slide master / slides) from the custom "Red Violet" / Integral
# palette over to the stock Office Theme palette, and switch the three
# tables that were still wearing the old custom table style over to the
# built-in "{37028D15-0C5D-4034-8BDC-C86BD1863D88}" table style.

$p = $ppt.ActivePresentation

# --- 1. Theme colours -------------------------------------------------
# ThemeColorScheme slot order is fixed: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1..accent6, 11 hlink, 12 folHlink. Setting it through any
# slide rewrites the single shared theme part, so slide 1 is as good as
# any other.
$officeThemeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388    # dk2      44546A
    4  = 15132391   # lt2      E7E6E6
    5  = 13998939   # accent1  5B9BD5
    6  = 3243501    # accent2  ED7D31
    7  = 10855845    # accent3  A5A5A5
    8  = 49407      # accent4  FFC000
    9  = 12874308   # accent5  4472C4
    10 = 4697456    # accent6  70AD47
    11 = 12673797   # hlink    0563C1
    12 = 7491477    # folHlink 954F72
}

$tcs = $p.Slides.Item(1).ThemeColorScheme
foreach ($slot in $officeThemeColors.Keys) {
    $tcs.Item($slot).RGB = $officeThemeColors[$slot]
}

# --- 2. Table styles ----------------------------------------------------
$newTableStyleId = "{37028D15-0C5D-4034-8BDC-C86BD1863D88}"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}
